# Added one more scenario for authentication (new row 6: Verify Authentication / TC_5)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new row 6 (new authentication test case) ---
$ws.Range("F6").Value = "InvalidPassword"
$ws.Range("B6").Value = "TC_5"

# --- Update existing rows: clarify the scenario descriptions ---
$ws.Range("C4").Value = "Verify when No Such Catalog in the site"
$ws.Range("C5").Value = "Verify when Execution is set as N"

# --- Continue filling in row 6 ---
$ws.Range("A6").Value = "Y"
$ws.Range("C6").Value = "Verify Authentication"
$ws.Range("D6").Value = "This will compare the product details from the listing page and details page"
$ws.Range("E6").Value = "lenproautomation8@lenqat.com"
$ws.Range("G6").Value = "Parts and Supplies"
$ws.Range("H6").Value = "Compressors:::Compressors"
$ws.Range("I6").Value = "Replace your compressor at LennoxPros.com."
$ws.Range("J6").Value = "10T46"
# "$1,173.00" looks like currency, force it to stay text like the other price cells
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "`$1,173.00"
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").Value = "H22J38BABC"

# Row height / wrap formatting to match the other data rows (A:K wrap, L plain)
$ws.Range("A6:K6").WrapText = $true
$ws.Rows("6:6").RowHeight = 45

# Hyperlink for the email address cell, styled like the other EMAILID cells
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:lenproautomation8@lenqat.com", "")
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E6").WrapText = $true

# Match the final selection left behind by the edit
$ws.Range("C6").Select()
